$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Rename header values: "CompanyMaster" -> "Company", "LocationMaster" -> "Location"
$ws.Range("AR1").Value = "Company"
$ws.Range("AS1").Value = "Location"

# Restore the original selection (AI1:BS1, active cell AI1) ...
$ws.Range("AI1:BS1").Select()
$excel.ActiveCell = $ws.Range("AI1")

# ... then scroll the window so the top-left visible column is X (was AC).
$excel.ActiveWindow.ScrollColumn = 24
$excel.ActiveWindow.ScrollRow = 1
